$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks 02-04 to 02-11")
$ws.Activate()

# Row 7 - "Create Base Enemy Class": fill in Time Spent / Over-Under
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 1

# Row 8 - "Create EnemyPrefab": fill in Time Spent / Over-Under, move status In Progress -> Done
$ws.Range("C8").Value = 0.5
$ws.Range("D8").Value = 0
$ws.Range("F2").Copy($ws.Range("F8"))

# Row 9 - "Create EnemySpawner": fill in Time Spent / Over-Under, move status In Progress -> Done
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 0
$ws.Range("F2").Copy($ws.Range("F9"))

# Row 12 - "Script Cleanup": update Time Estimated / Time Spent
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 1.5

# Update the last active selection to D4
$ws.Range("D4").Select() | Out-Null
